# Daily attendance processing - 2025-11-25 07:48:54
# Reorders the comma-separated "Recorded By" names in column G for the
# affected rows, per the source export's updated ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2 = @{ Old = "system, backup@backdoor.com, System"; New = "backup@backdoor.com, System, system" }
    3 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    6 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    10 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    11 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    12 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    13 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    14 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    15 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    17 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    18 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    19 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    20 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    21 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    22 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    24 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    26 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    28 = @{ Old = "system, backup@backdoor.com, System"; New = "backup@backdoor.com, System, system" }
    29 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    32 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    36 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    37 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    38 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    39 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    40 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    41 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    43 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    44 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    45 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    46 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    47 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    48 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    50 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    52 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    54 = @{ Old = "system, backup@backdoor.com, System"; New = "backup@backdoor.com, System, system" }
    55 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    58 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    62 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    63 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    64 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    65 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    66 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    67 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    69 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    70 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    71 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    72 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    73 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    74 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    76 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    78 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    83 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    84 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    85 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    86 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    87 = @{ Old = "dnasr281@gmail.com, admin@admin.com"; New = "admin@admin.com, dnasr281@gmail.com" }
    90 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    92 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    93 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    94 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    96 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    99 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    101 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    109 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    110 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    111 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    112 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    113 = @{ Old = "dnasr281@gmail.com, admin@admin.com"; New = "admin@admin.com, dnasr281@gmail.com" }
    116 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    118 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    119 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    120 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    122 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    125 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    127 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    135 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    136 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    137 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    138 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    139 = @{ Old = "dnasr281@gmail.com, admin@admin.com"; New = "admin@admin.com, dnasr281@gmail.com" }
    142 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    144 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    145 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    146 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    148 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    151 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
    153 = @{ Old = "dnasr281@gmail.com, System"; New = "System, dnasr281@gmail.com" }
}

$applied = 0
$skipped = 0

foreach ($row in $changes.Keys) {
    $entry = $changes[$row]
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value()
    if ($current -eq $entry.Old) {
        $cell.Value = $entry.New
        $applied++
    } else {
        Write-Host "Row $row : unexpected value [$current], expected [$($entry.Old)] - skipped"
        $skipped++
    }
}

Write-Host "Applied $applied changes, skipped $skipped"
